# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Macroferia Regional de Talca - Zanahoria"
# before the existing row 210, shifting the subsequent rows down (sheet grows
# from A1:R222 to A1:R223).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row at position 210; everything from the old row 210
# downward (through 222) shifts down to 211..223.
$ws.Rows.Item(210).Insert()

# Populate the newly inserted row 210 with the new week's data.
$ws.Range("A210").Value = 5
$ws.Range("B210").Value = "Macroferia Regional de Talca"
$ws.Range("C210").Value = "Maule"
$ws.Range("D210").Value = 44516
$ws.Range("E210").Value = 7
$ws.Range("F210").Value = 100114013
$ws.Range("G210").Value = "Zanahoria"
$ws.Range("H210").Value = "Sin especificar"
$ws.Range("I210").Value = "Primera"
$ws.Range("J210").Value = 400
$ws.Range("K210").Value = 9000
$ws.Range("L210").Value = 9000
$ws.Range("M210").Value = 9000
$ws.Range("N210").Value = "$/saco 20 kilos"
$ws.Range("O210").Value = "Provincia del Elquí"
$ws.Range("P210").Value = 450
$ws.Range("Q210").Value = 20
$ws.Range("R210").Value = "Hortaliza"
